# LevelInfo.xlsx: remove the "升级奇术" (arcane/technique upgrade) level-up
# tip row from the LevelInfo table. Deleting the whole sheet row shifts the
# rows below it up by one and keeps the table ("表1") / autofilter range in
# sync automatically (Excel shrinks ref="A1:E13" -> "A1:E12").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Cells.Find("升级奇术")
if ($target -ne $null) {
    $target.EntireRow.Delete() | Out-Null
}

# Matches the author's resulting selection in the saved file.
$ws.Range("D8").Select() | Out-Null
